$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row at position 25
$ws.Rows.Item(25).Insert()
$ws.Range("A25").Value = 3
$ws.Range("B25").Value = 'Femacal de La Calera'
$ws.Range("C25").Value = 'Coquimbo'
$ws.Range("D25").Value = 44425
$ws.Range("E25").Value = 5
$ws.Range("F25").Value = 100112039
$ws.Range("G25").Value = 'Ciboulette'
$ws.Range("H25").Value = 'Sin especificar'
$ws.Range("I25").Value = 'Primera'
$ws.Range("J25").Value = 160
$ws.Range("K25").Value = 1500
$ws.Range("L25").Value = 1500
$ws.Range("M25").Value = 1500
$ws.Range("N25").Value = '$/docena de atados'
$ws.Range("O25").Value = 'Provincia de Quillota'
$ws.Range("P25").Value = 500
$ws.Range("Q25").Value = 3
$ws.Range("R25").Value = 'Hortaliza'

# Insert new row at position 67
$ws.Rows.Item(67).Insert()
$ws.Range("A67").Value = 3
$ws.Range("B67").Value = 'Femacal de La Calera'
$ws.Range("C67").Value = 'Coquimbo'
$ws.Range("D67").Value = 44427
$ws.Range("E67").Value = 5
$ws.Range("F67").Value = 100112039
$ws.Range("G67").Value = 'Ciboulette'
$ws.Range("H67").Value = 'Sin especificar'
$ws.Range("I67").Value = 'Primera'
$ws.Range("J67").Value = 160
$ws.Range("K67").Value = 1500
$ws.Range("L67").Value = 1500
$ws.Range("M67").Value = 1500
$ws.Range("N67").Value = '$/docena de atados'
$ws.Range("O67").Value = 'Provincia de Quillota'
$ws.Range("P67").Value = 500
$ws.Range("Q67").Value = 3
$ws.Range("R67").Value = 'Hortaliza'

# Insert new row at position 70
$ws.Rows.Item(70).Insert()
$ws.Range("A70").Value = 3
$ws.Range("B70").Value = 'Femacal de La Calera'
$ws.Range("C70").Value = 'Coquimbo'
$ws.Range("D70").Value = 44421
$ws.Range("E70").Value = 5
$ws.Range("F70").Value = 100112039
$ws.Range("G70").Value = 'Ciboulette'
$ws.Range("H70").Value = 'Sin especificar'
$ws.Range("I70").Value = 'Primera'
$ws.Range("J70").Value = 180
$ws.Range("K70").Value = 1500
$ws.Range("L70").Value = 1500
$ws.Range("M70").Value = 1500
$ws.Range("N70").Value = '$/docena de atados'
$ws.Range("O70").Value = 'Provincia de Quillota'
$ws.Range("P70").Value = 500
$ws.Range("Q70").Value = 3
$ws.Range("R70").Value = 'Hortaliza'

# Insert new row at position 76
$ws.Rows.Item(76).Insert()
$ws.Range("A76").Value = 3
$ws.Range("B76").Value = 'Femacal de La Calera'
$ws.Range("C76").Value = 'Coquimbo'
$ws.Range("D76").Value = 44426
$ws.Range("E76").Value = 5
$ws.Range("F76").Value = 100112039
$ws.Range("G76").Value = 'Ciboulette'
$ws.Range("H76").Value = 'Sin especificar'
$ws.Range("I76").Value = 'Primera'
$ws.Range("J76").Value = 160
$ws.Range("K76").Value = 1500
$ws.Range("L76").Value = 1500
$ws.Range("M76").Value = 1500
$ws.Range("N76").Value = '$/docena de atados'
$ws.Range("O76").Value = 'Provincia de Quillota'
$ws.Range("P76").Value = 500
$ws.Range("Q76").Value = 3
$ws.Range("R76").Value = 'Hortaliza'

# Insert new row at position 113
$ws.Rows.Item(113).Insert()
$ws.Range("A113").Value = 3
$ws.Range("B113").Value = 'Femacal de La Calera'
$ws.Range("C113").Value = 'Coquimbo'
$ws.Range("D113").Value = 44417
$ws.Range("E113").Value = 5
$ws.Range("F113").Value = 100112039
$ws.Range("G113").Value = 'Ciboulette'
$ws.Range("H113").Value = 'Sin especificar'
$ws.Range("I113").Value = 'Primera'
$ws.Range("J113").Value = 160
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1500
$ws.Range("N113").Value = '$/docena de atados'
$ws.Range("O113").Value = 'Provincia de Quillota'
$ws.Range("P113").Value = 500
$ws.Range("Q113").Value = 3
$ws.Range("R113").Value = 'Hortaliza'

# Insert new row at position 114
$ws.Rows.Item(114).Insert()
$ws.Range("A114").Value = 3
$ws.Range("B114").Value = 'Femacal de La Calera'
$ws.Range("C114").Value = 'Coquimbo'
$ws.Range("D114").Value = 44419
$ws.Range("E114").Value = 5
$ws.Range("F114").Value = 100112039
$ws.Range("G114").Value = 'Ciboulette'
$ws.Range("H114").Value = 'Sin especificar'
$ws.Range("I114").Value = 'Primera'
$ws.Range("J114").Value = 130
$ws.Range("K114").Value = 1500
$ws.Range("L114").Value = 1500
$ws.Range("M114").Value = 1500
$ws.Range("N114").Value = '$/docena de atados'
$ws.Range("O114").Value = 'Provincia de Quillota'
$ws.Range("P114").Value = 500
$ws.Range("Q114").Value = 3
$ws.Range("R114").Value = 'Hortaliza'

# Insert new row at position 115
$ws.Rows.Item(115).Insert()
$ws.Range("A115").Value = 3
$ws.Range("B115").Value = 'Femacal de La Calera'
$ws.Range("C115").Value = 'Coquimbo'
$ws.Range("D115").Value = 44420
$ws.Range("E115").Value = 5
$ws.Range("F115").Value = 100112039
$ws.Range("G115").Value = 'Ciboulette'
$ws.Range("H115").Value = 'Sin especificar'
$ws.Range("I115").Value = 'Primera'
$ws.Range("J115").Value = 160
$ws.Range("K115").Value = 1500
$ws.Range("L115").Value = 1500
$ws.Range("M115").Value = 1500
$ws.Range("N115").Value = '$/docena de atados'
$ws.Range("O115").Value = 'Provincia de Quillota'
$ws.Range("P115").Value = 500
$ws.Range("Q115").Value = 3
$ws.Range("R115").Value = 'Hortaliza'

# Insert new row at position 157
$ws.Rows.Item(157).Insert()
$ws.Range("A157").Value = 3
$ws.Range("B157").Value = 'Femacal de La Calera'
$ws.Range("C157").Value = 'Coquimbo'
$ws.Range("D157").Value = 44418
$ws.Range("E157").Value = 5
$ws.Range("F157").Value = 100112039
$ws.Range("G157").Value = 'Ciboulette'
$ws.Range("H157").Value = 'Sin especificar'
$ws.Range("I157").Value = 'Primera'
$ws.Range("J157").Value = 150
$ws.Range("K157").Value = 1500
$ws.Range("L157").Value = 1500
$ws.Range("M157").Value = 1500
$ws.Range("N157").Value = '$/docena de atados'
$ws.Range("O157").Value = 'Provincia de Quillota'
$ws.Range("P157").Value = 500
$ws.Range("Q157").Value = 3
$ws.Range("R157").Value = 'Hortaliza'
